# Applies the v1.2.1 -> v1.2.3 changes to "UC005 - Listar Empenhos Pendentes--GTP-.xlsx"
#
# 1) TC2's expected result (D18) gains a clause about ordering by "numero da diaria".
# 2) TC3 and TC4 swap their step content:
#    - TC3 (rows 22-26) used to be about atribuir/desatribuir; it now becomes
#      "realizar o empenho de uma diária" / "Apresenta a tela de Registrar Empenho".
#    - TC4 (rows 29-33) used to be about realizar o empenho; it now becomes
#      "atribuir/desatribuir" / "Atualiza a lista de registros de solicitações...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update TC2's "Expected Results" cell (D18)
$ws.Range("D18").Value = "SYSTEM Exibe a lista de solicitações aguardando serem empenhadas ordenado pelo numero da diaria em ordem crescente."

# 2) Swap the step content between TC3 (row 26) and TC4 (row 33)
$ws.Range("B26").Value = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
$ws.Range("D26").Value = "SYSTEM Apresenta a tela de Registrar Empenho"

$ws.Range("B33").Value = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D33").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
